# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" everywhere it is used
#    (Overview!E:F, zh-cn!C, de-de!C).
# 2. Narrow the "Status" column(s) from ~17.22 chars to ~13.41 chars:
#    Overview columns E & F, and column C on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Replace the status text wherever it occurs -------------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $cellValue = $cell.Value2
        if ("Ready for handoff" -eq $cellValue) {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Resize the Status columns -------------------------------------------
# The target raw width (13.4101845877511 "XML" units) is not exactly reachable
# through the ColumnWidth COM property (it snaps to the nearest 1/6th), so we
# use the closest input that lands on the nearest achievable width.
$targetColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth   # column F (de-de status)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth       # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth       # column C (Status)
